$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("windspeed")

# Update column B values for rows 13 through 32 to 15
$ws.Range("B13:B32").Value = 15

# Update the selected range/active cell on the sheet
$ws.Range("B23:B32").Select()

